# Update countries & provincias Spain
#
# 1) Swap the country labels for rows 37/38 (Corea del Sur <-> Emiratos Arabes
#    Unidos changed position in the shared country list), while each row
#    keeps its own (updated) statistics.
# 2) Refresh the numeric statistics (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#    affected countries: Estados Unidos, Reino Unido, Canada, Corea del Sur /
#    Emiratos Arabes Unidos, Irak and Sri Lanka.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 / 38: swap the two country names, and give each row its own
#     refreshed numbers ------------------------------------------------
$ws.Range("A37").Value = "Emiratos Arabes Unidos"
$ws.Range("A38").Value = "Corea del Sur"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 999390
$ws.Range("C4").Value = 12230
$ws.Range("D4").Value = 137271
$ws.Range("E4").Value = 805944
$ws.Range("F4").Value = 14175
$ws.Range("G4").Value = 762
$ws.Range("H4").Value = 56175

# Row 9 - Reino Unido
$ws.Range("B9").Value = 157149
$ws.Range("C9").Value = 4309
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 135713
$ws.Range("F9").Value = 1559
$ws.Range("G9").Value = 360
$ws.Range("H9").Value = 21092

# Row 15 - Canada
$ws.Range("B15").Value = 48229
$ws.Range("C15").Value = 1334
$ws.Range("D15").Value = 17916
$ws.Range("E15").Value = 27612
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 141
$ws.Range("H15").Value = 2701

# Row 37 - Emiratos Arabes Unidos (new numbers)
$ws.Range("B37").Value = 10839
$ws.Range("C37").Value = 490
$ws.Range("D37").Value = 2090
$ws.Range("E37").Value = 8667
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 82

# Row 38 - Corea del Sur (numbers formerly held by row 37)
$ws.Range("B38").Value = 10738
$ws.Range("C38").Value = 10
$ws.Range("D38").Value = 8764
$ws.Range("E38").Value = 1731
$ws.Range("F38").Value = 55
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 243

# Row 69 - Irak
$ws.Range("B69").Value = 1847
$ws.Range("C69").Value = 27
$ws.Range("D69").Value = 1286
$ws.Range("E69").Value = 473
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 88

# Row 105 - Sri Lanka (provincias Spain block in commit message, but the
# actual changed row here is the countries sheet row 105)
$ws.Range("B105").Value = 588
$ws.Range("C105").Value = 65
$ws.Range("D105").Value = 126
$ws.Range("E105").Value = 455
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7
